# Updated cryptos list on Fri Aug  4 08:46:18 UTC 2023 with GitHub Actions
# Refresh the scraped Price (D) and Volume(1h) (E) columns, and fix the
# ordering of the Aave / BabyDogeCoin rows (48-49 swapped positions).
#
# Values are written with a leading apostrophe to force Excel to treat
# the numeric-looking price strings (e.g. "241.95", "1.000",
# "0.00000000120") as literal text instead of auto-converting them to
# numbers, matching how the source data is stored (plain text cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.179.04"
$ws.Range("E2").Value = "'  +0.43%  "
$ws.Range("D3").Value = "'1.835.61"
$ws.Range("E3").Value = "'  +0.34%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'241.95"
$ws.Range("E5").Value = "'  +1.27%  "
$ws.Range("D6").Value = "'0.6581"
$ws.Range("E6").Value = "'  -0.44%  "
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("D8").Value = "'0.07426"
$ws.Range("E8").Value = "'  +1.35%  "
$ws.Range("D9").Value = "'0.2927"
$ws.Range("E9").Value = "'  -0.52%  "
$ws.Range("D10").Value = "'22.95"
$ws.Range("E10").Value = "'  +1.28%  "
$ws.Range("D11").Value = "'0.07773"
$ws.Range("E11").Value = "'  +1.79%  "
$ws.Range("D12").Value = "'1.863.24"
$ws.Range("E12").Value = "'  +1.63%  "
$ws.Range("E13").Value = "'  -0.45%  "
$ws.Range("D14").Value = "'0.6653"
$ws.Range("E14").Value = "'  -0.93%  "
$ws.Range("D15").Value = "'82.82"
$ws.Range("E15").Value = "'  -3.80%  "
$ws.Range("D16").Value = "'6.112"
$ws.Range("E16").Value = "'  +0.12%  "
$ws.Range("D17").Value = "'0.000008599"
$ws.Range("E17").Value = "'  +4.72%  "
$ws.Range("D18").Value = "'29.203.15"
$ws.Range("E18").Value = "'  +0.54%  "
$ws.Range("D19").Value = "'2.149.88"
$ws.Range("E19").Value = "'  +3.64%  "
$ws.Range("D20").Value = "'226.63"
$ws.Range("E20").Value = "'  -0.40%  "
$ws.Range("D21").Value = "'12.45"
$ws.Range("E21").Value = "'  +0.12%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "'  +0.15%  "
$ws.Range("D23").Value = "'7.121"
$ws.Range("E23").Value = "'  -1.71%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "'  +0.05%  "
$ws.Range("D25").Value = "'159.31"
$ws.Range("E25").Value = "'  -0.71%  "
$ws.Range("D26").Value = "'8.603"
$ws.Range("E26").Value = "'  -0.60%  "
$ws.Range("D27").Value = "'0.1395"
$ws.Range("E27").Value = "'  -1.78%  "
$ws.Range("D28").Value = "'17.93"
$ws.Range("E28").Value = "'  +0.04%  "
$ws.Range("D29").Value = "'1.514"
$ws.Range("E29").Value = "'  +0.99%  "
$ws.Range("D30").Value = "'4.114"
$ws.Range("E30").Value = "'  -2.41%  "
$ws.Range("D31").Value = "'4.044"
$ws.Range("E31").Value = "'  -1.32%  "
$ws.Range("E32").Value = "'  -0.22%  "
$ws.Range("D33").Value = "'0.05269"
$ws.Range("E33").Value = "'  -0.61%  "
$ws.Range("D34").Value = "'1.866"
$ws.Range("E34").Value = "'  +1.12%  "
$ws.Range("D35").Value = "'0.7397"
$ws.Range("E35").Value = "'  -1.25%  "
$ws.Range("E36").Value = "'  +1.87%  "
$ws.Range("D38").Value = "'1.305.42"
$ws.Range("E38").Value = "'  +1.04%  "
$ws.Range("D39").Value = "'0.01797"
$ws.Range("E39").Value = "'  -0.50%  "
$ws.Range("E40").Value = "'  +1.07%  "
$ws.Range("D41").Value = "'0.9209"
$ws.Range("E41").Value = "'  -0.15%  "
$ws.Range("D42").Value = "'6.047"
$ws.Range("E42").Value = "'  +1.05%  "
$ws.Range("D43").Value = "'0.08636"
$ws.Range("E43").Value = "'  +13.35%  "
$ws.Range("E44").Value = "'  +0.12%  "
$ws.Range("D45").Value = "'102.46"
$ws.Range("E45").Value = "'  -0.99%  "
$ws.Range("D46").Value = "'2.027.60"
$ws.Range("E46").Value = "'  +2.67%  "
$ws.Range("D47").Value = "'0.5144"
$ws.Range("E47").Value = "'  -0.62%  "
$ws.Range("B48").Value = "'BabyDogeCoin"
$ws.Range("C48").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000120"
$ws.Range("E48").Value = "'  -2.58%  "
$ws.Range("B49").Value = "'Aave"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'63.57"
$ws.Range("E49").Value = "'  +0.50%  "
$ws.Range("D50").Value = "'1.751"
$ws.Range("D51").Value = "'0.05847"
$ws.Range("E51").Value = "'  -1.07%  "
